$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a new column before column C, shifting old C..Q to D..R
$ws.Columns("C").Insert()

$ws.Range("B13").Value = "S011"
$ws.Range("C3").Value = "Recency"
$ws.Range("C13").Value = "MRE"
$ws.Range("D13").Value = "Bogo's Club"
$ws.Range("E13").Value = 41296
$ws.Range("F13").Value = 40817
$ws.Range("G13").Value = 41182

$ws.Range("D1").ColumnWidth = 13.21875
$ws.Range("P1").ColumnWidth = 9.5546875

# Comments don't auto-shift with the column insert, so relocate them manually
# (move right-most first to avoid collisions)
$txtF = $ws.Range("F3").Comment.Text()
$ws.Range("F3").Comment.Delete()
$ws.Range("G3").AddComment($txtF)

$txtE = $ws.Range("E3").Comment.Text()
$ws.Range("E3").Comment.Delete()
$ws.Range("F3").AddComment($txtE)

$txtD = $ws.Range("D3").Comment.Text()
$ws.Range("D3").Comment.Delete()
$ws.Range("E3").AddComment($txtD)

$ws.Range("B8").Select()
